$d = $word.ActiveDocument

# Only the "2022 observation dates" heading runs (which read "...για τον
# αστερισμό του Αστερισμός Ταύρου...") should lose the redundant
# "αστερισμό του " wording; similar phrasing about the constellation Orion
# elsewhere in the document must stay untouched. Matching on the unique
# surrounding text "τον αστερισμό του Αστερισμός" -> "τον  Αστερισμός"
# (note the double space is intentional, matching the target text) keeps
# the edit scoped to exactly those 4 heading runs.
$d.Content.Find.Execute("τον αστερισμό του Αστερισμός", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "τον  Αστερισμός", 2)
